# Update automatico via Actualizar 05-12-2020 02-47-07
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HOSPITALES")

# Row 152 (Cortés department municipalities)
$ws.Range("I152").Value = 3
$ws.Range("J152").Value = '0503'
$ws.Range("K152").Value = 'Omoa'
$ws.Range("M152").Value = 1
$ws.Range("N152").Value = '050301'
$ws.Range("O152").Value = 'Omoa'
$ws.Range("Q152").Value = 'HND-0503'
$ws.Range("V152").Value = 15.774128
$ws.Range("W152").Value = -88.038326

# Row 153 (Cortés department municipalities)
$ws.Range("I153").Value = 3
$ws.Range("J153").Value = '0503'
$ws.Range("K153").Value = 'Omoa'
$ws.Range("M153").Value = 7
$ws.Range("N153").Value = '050307'
$ws.Range("O153").Value = 'Cuyamel'
$ws.Range("Q153").Value = 'HND-0503'
$ws.Range("V153").Value = 15.663361
$ws.Range("W153").Value = -88.194594

# Row 154 (Cortés department municipalities)
$ws.Range("I154").Value = 2
$ws.Range("J154").Value = '0502'
$ws.Range("K154").Value = 'Choloma'
$ws.Range("M154").Value = 1
$ws.Range("N154").Value = '050201'
$ws.Range("O154").Value = 'Col. Lopez Arellano'
$ws.Range("Q154").Value = 'HND-0502'
$ws.Range("V154").Value = 15.610602
$ws.Range("W154").Value = -87.951762

# Row 155 (Cortés department municipalities)
$ws.Range("I155").Value = 12
$ws.Range("J155").Value = '0512'
$ws.Range("K155").Value = 'La Lima'
$ws.Range("M155").Value = 1
$ws.Range("N155").Value = '051201'
$ws.Range("O155").Value = 'La Lima'
$ws.Range("Q155").Value = 'HND-0512'
$ws.Range("V155").Value = 15.439396
$ws.Range("W155").Value = -87.928896

# Row 156 (Cortés department municipalities)
$ws.Range("I156").Value = 9
$ws.Range("J156").Value = '0509'
$ws.Range("K156").Value = 'San Manuel'
$ws.Range("M156").Value = 1
$ws.Range("N156").Value = '050901'
$ws.Range("O156").Value = 'San Manuel'
$ws.Range("Q156").Value = 'HND-0509'
$ws.Range("V156").Value = 15.329429
$ws.Range("W156").Value = -87.921024

# Row 157 (Cortés department municipalities)
$ws.Range("I157").Value = 11
$ws.Range("J157").Value = '0511'
$ws.Range("K157").Value = 'Villanueva'
$ws.Range("M157").Value = 1
$ws.Range("N157").Value = '051101'
$ws.Range("O157").Value = 'Villanueva'
$ws.Range("Q157").Value = 'HND-0511'
$ws.Range("V157").Value = 15.312935
$ws.Range("W157").Value = -87.993704

# Row 158 (Cortés department municipalities)
$ws.Range("I158").Value = 5
$ws.Range("J158").Value = '0505'
$ws.Range("K158").Value = 'Potrerillos'
$ws.Range("M158").Value = 1
$ws.Range("N158").Value = '050501'
$ws.Range("O158").Value = 'Potrerillos'
$ws.Range("Q158").Value = 'HND-0505'
$ws.Range("V158").Value = 15.228071
$ws.Range("W158").Value = -87.964017

# Row 159 (Cortés department municipalities)
$ws.Range("I159").Value = 7
$ws.Range("J159").Value = '0507'
$ws.Range("K159").Value = 'San Antonio de Cortés'
$ws.Range("M159").Value = 1
$ws.Range("N159").Value = '050701'
$ws.Range("O159").Value = 'San Antonio de Cortés'
$ws.Range("Q159").Value = 'HND-0507'
$ws.Range("V159").Value = 15.114108
$ws.Range("W159").Value = -88.040539

# Row 160 (Cortés department municipalities)
$ws.Range("E160").Value = 5
$ws.Range("F160").Value = '05'
$ws.Range("G160").Value = 'Cortés'
$ws.Range("I160").Value = 10
$ws.Range("J160").Value = '0510'
$ws.Range("K160").Value = 'Santa Cruz de Yojoa'
$ws.Range("M160").Value = 1
$ws.Range("N160").Value = '051001'
$ws.Range("O160").Value = 'Santa Cruz de Yojoa'
$ws.Range("Q160").Value = 'HND-0510'
$ws.Range("V160").Value = 14.97907
$ws.Range("W160").Value = -87.890919

# Row 161 (Cortés department municipalities)
$ws.Range("Q161").Value = 'HND-'

# N161 gets touched (empty text-formatted cell) without content, matching source edit
$ws.Range("N161").NumberFormat = "@"

# Update view/selection to reflect where the edit left the cursor
$w = $excel.Windows.Item(1)
$w.ScrollRow = 154
$w.ScrollColumn = 21
$ws.Range("W162").Select()
